$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Give the new header cells the same formatting as the existing header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data cells I2 and J2 (plain, unstyled, numeric - like the rest of row 2)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4
